$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lines = @(
    "填写须知：",
    "1、不能增加、删除列；",
    "2、不能修改灰色单元格；",
    "3、红色字段为必填字段，黑色字段为选填字段；",
    "4、删除餐厅编码后再导入，系统会根据ID删除没有餐厅编码的记录；",
    "5、删除餐厅编码后再导入，系统会根据ID删除没有餐厅编码的记录；",
    "6、删除餐厅编码后再导入，系统会根据ID删除没有餐厅编码的记录；删除餐厅编码后再导入，系统会根据ID删除没有餐厅编码的记录删除餐厅编码后再导入，系统会根据ID删除没有餐厅编码的记录",
    "7、删除餐厅编码后再导入，系统会根据ID删除没有餐厅编码的记录；"
)
$text = [string]::Join([char]10, $lines)

$ws.Range("A1").Value = $text

$ws.Rows.Item(1).RowHeight = 192.0
